# Update "想去人数" (column F) figures across all sheets to match the
# freshly generated data output (gh-pages build at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 1312
$ws1.Range("F6").Value  = 254
$ws1.Range("F7").Value  = 395
$ws1.Range("F8").Value  = 8389
$ws1.Range("F10").Value = 10349
$ws1.Range("F11").Value = 90
$ws1.Range("F23").Value = 401
$ws1.Range("F25").Value = 1765
$ws1.Range("F27").Value = 520
$ws1.Range("F29").Value = 279
$ws1.Range("F33").Value = 1103
$ws1.Range("F37").Value = 431
$ws1.Range("F39").Value = 278
$ws1.Range("F42").Value = 502
$ws1.Range("F43").Value = 333
$ws1.Range("F44").Value = 80
$ws1.Range("F45").Value = 270
$ws1.Range("F48").Value = 69
$ws1.Range("F49").Value = 72

# Sheet "演出" (sheetId 2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 37

# Sheet "本地生活" (sheetId 3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 205

# Sheet "全部类型" (sheetId 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 205
$ws4.Range("F7").Value  = 1312
$ws4.Range("F9").Value  = 254
$ws4.Range("F12").Value = 395
$ws4.Range("F13").Value = 8389
$ws4.Range("F15").Value = 10349
$ws4.Range("F16").Value = 90
$ws4.Range("F23").Value = 1765
$ws4.Range("F24").Value = 520
$ws4.Range("F26").Value = 279
$ws4.Range("F29").Value = 37
$ws4.Range("F30").Value = 1103
$ws4.Range("F35").Value = 431
$ws4.Range("F40").Value = 502
$ws4.Range("F42").Value = 333
$ws4.Range("F43").Value = 80
$ws4.Range("F44").Value = 270
